$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 102.71429
$ws.Range("I11").Value = 102.71429
$ws.Range("K11").Value = 102.71429
$ws.Range("M11").Value = 37.28570999999999
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H31").Value = 744.8
$ws.Range("I31").Value = 744.8
$ws.Range("K31").Value = 2234.4
$ws.Range("M31").Value = -2004.4
$ws.Range("H38").Value = 1434
$ws.Range("I38").Value = 23
$ws.Range("J38").Value = 9900
$ws.Range("K38").Value = 69
$ws.Range("L38").Value = 29700
$ws.Range("M38").Value = 303
$ws.Range("N38").Value = -30444
$ws.Range("H39").Value = 148.375
$ws.Range("I39").Value = 148.375
$ws.Range("K39").Value = 445.125
$ws.Range("M39").Value = -149.125
$ws.Range("H86").Value = 3500.5
$ws.Range("I86").Value = 3500.5
$ws.Range("K86").Value = 3500.5
$ws.Range("M86").Value = -2377.5
$ws.Range("H89").Value = 3500.5
$ws.Range("I89").Value = 3500.5
$ws.Range("K89").Value = 17502.5
$ws.Range("M89").Value = -11886.5
$ws.Range("H112").Value = 4500
$ws.Range("J112").Value = 4500
$ws.Range("L112").Value = 13500
$ws.Range("N112").Value = -15716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2127.5715
$ws.Range("I2").Value = 2148.8333
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2148.8333
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -2035.8333
$ws.Range("N2").Value = -2226
$ws.Range("H13").Value = 2000
$ws.Range("I13").Value = 2000
$ws.Range("K13").Value = 2000
$ws.Range("M13").Value = -1856
$ws.Range("H32").Value = 1149.5
$ws.Range("I32").Value = 1149.5
$ws.Range("K32").Value = 1149.5
$ws.Range("M32").Value = -862.5
$ws.Range("H63").Value = 6487.5
$ws.Range("J63").Value = 500
$ws.Range("L63").Value = 500
$ws.Range("N63").Value = -1872
$ws.Range("H66").Value = 6487.5
$ws.Range("J66").Value = 500
$ws.Range("L66").Value = 2500
$ws.Range("N66").Value = -9364
$ws.Range("H88").Value = 718.5
$ws.Range("I88").Value = 800
$ws.Range("J88").Value = 691.3333
$ws.Range("K88").Value = 800
$ws.Range("L88").Value = 691.3333
$ws.Range("M88").Value = -394
$ws.Range("N88").Value = -1503.3333
$ws.Range("H91").Value = 718.5
$ws.Range("I91").Value = 800
$ws.Range("J91").Value = 691.3333
$ws.Range("K91").Value = 800
$ws.Range("L91").Value = 691.3333
$ws.Range("M91").Value = 604
$ws.Range("N91").Value = -3499.3333
$ws.Range("H116").Value = 2127.5715
$ws.Range("I116").Value = 2148.8333
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2148.8333
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 145.1667000000002
$ws.Range("N116").Value = -6588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2127.5715
$ws.Range("I3").Value = 2148.8333
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2148.8333
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -2034.8333
$ws.Range("N3").Value = -2228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 100
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 100
$ws.Range("M12").Value = 70
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 798.4
$ws.Range("I8").Value = 798.4
$ws.Range("K8").Value = 2395.2
$ws.Range("M8").Value = -2256.2
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 1
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 170

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4997.5
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 15337.333
$ws.Range("I132").Value = 13012
$ws.Range("K132").Value = 39036
$ws.Range("M132").Value = -36506

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -7876
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
$ws.Range("H94").Value = 59999
$ws.Range("J94").Value = 59999
$ws.Range("L94").Value = 59999
$ws.Range("N94").Value = -61351
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 10000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 30000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1999.6666
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 2997.5
$ws.Range("K17").Value = 4
$ws.Range("L17").Value = 2997.5
$ws.Range("M17").Value = 168
$ws.Range("N17").Value = -3341.5
$ws.Range("H38").Value = 1500
$ws.Range("I38").Value = 1500
$ws.Range("K38").Value = 1500
$ws.Range("M38").Value = -1027
$ws.Range("H69").Value = 26249.25
$ws.Range("J69").Value = 26249.25
$ws.Range("L69").Value = 26249.25
$ws.Range("N69").Value = -27747.25
$ws.Range("H72").Value = 26249.25
$ws.Range("J72").Value = 26249.25
$ws.Range("L72").Value = 78747.75
$ws.Range("N72").Value = -86235.75
